$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the string value in C1 ("asdasda") with the numeric value 201005
$ws.Range("C1").Value = 201005
